$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 15, pushing the existing row 15 down to row 17
$ws.Rows.Item(15).Resize(2).Insert()

# New row 15 values
$ws.Cells.Item(15, 1).Value = 8
$ws.Cells.Item(15, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(15, 3).Value = "Coquimbo"
$ws.Cells.Item(15, 4).Value = 44610
$ws.Cells.Item(15, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(15, 5).Value = 4
$ws.Cells.Item(15, 6).Value = "Fruta"
$ws.Cells.Item(15, 7).Value = 100107
$ws.Cells.Item(15, 8).Value = "Otros"
$ws.Cells.Item(15, 9).Value = 100107011
$ws.Cells.Item(15, 10).Value = "Tuna"
$ws.Cells.Item(15, 11).Value = "Sin especificar"
$ws.Cells.Item(15, 12).Value = "Primera"
$ws.Cells.Item(15, 13).Value = 200
$ws.Cells.Item(15, 14).Value = 13000
$ws.Cells.Item(15, 15).Value = 14000
$ws.Cells.Item(15, 16).Value = 13500
$ws.Cells.Item(15, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(15, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(15, 19).Value = 750
$ws.Cells.Item(15, 20).Value = 18

# New row 16 values
$ws.Cells.Item(16, 1).Value = 8
$ws.Cells.Item(16, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(16, 3).Value = "Coquimbo"
$ws.Cells.Item(16, 4).Value = 44610
$ws.Cells.Item(16, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(16, 5).Value = 4
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100107
$ws.Cells.Item(16, 8).Value = "Otros"
$ws.Cells.Item(16, 9).Value = 100107011
$ws.Cells.Item(16, 10).Value = "Tuna"
$ws.Cells.Item(16, 11).Value = "Sin especificar"
$ws.Cells.Item(16, 12).Value = "Segunda"
$ws.Cells.Item(16, 13).Value = 200
$ws.Cells.Item(16, 14).Value = 11000
$ws.Cells.Item(16, 15).Value = 12000
$ws.Cells.Item(16, 16).Value = 11500
$ws.Cells.Item(16, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(16, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(16, 19).Value = 639
$ws.Cells.Item(16, 20).Value = 18
